$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Simple Price / Volume updates ---
Set-TextCell $ws.Range("D2") "51.872.75"
Set-TextCell $ws.Range("E2") "  -0.25%  "
Set-TextCell $ws.Range("D3") "2.828.91"
Set-TextCell $ws.Range("E3") "  +1.76%  "
Set-TextCell $ws.Range("E4") "  +0.06%  "
Set-TextCell $ws.Range("D5") "357.42"
Set-TextCell $ws.Range("E5") "  +4.33%  "
Set-TextCell $ws.Range("D6") "112.36"
Set-TextCell $ws.Range("E6") "  -2.46%  "
Set-TextCell $ws.Range("E7") "  +3.82%  "
Set-TextCell $ws.Range("E9") "  +4.18%  "
Set-TextCell $ws.Range("D10") "41.13"
Set-TextCell $ws.Range("E10") "  -1.73%  "
Set-TextCell $ws.Range("D11") "0.0858"
Set-TextCell $ws.Range("E11") "  +0.14%  "
Set-TextCell $ws.Range("D12") "20.10"
Set-TextCell $ws.Range("E12") "  +0.38%  "
Set-TextCell $ws.Range("E13") "  +1.14%  "
Set-TextCell $ws.Range("D14") "7.78"
Set-TextCell $ws.Range("E14") "  +1.99%  "
Set-TextCell $ws.Range("D15") "3.278.07"
Set-TextCell $ws.Range("E15") "  +1.96%  "
Set-TextCell $ws.Range("D16") "2.827.78"
Set-TextCell $ws.Range("E17") "  +6.00%  "
Set-TextCell $ws.Range("D18") "51.802.99"
Set-TextCell $ws.Range("E18") "  -0.12%  "
Set-TextCell $ws.Range("D19") "7.53"
Set-TextCell $ws.Range("E19") "  +7.20%  "
Set-TextCell $ws.Range("D20") "3.15"
Set-TextCell $ws.Range("E20") "  -1.21%  "
Set-TextCell $ws.Range("D21") "13.43"
Set-TextCell $ws.Range("E21") "  +1.46%  "
Set-TextCell $ws.Range("D22") "0.0₃0991"
Set-TextCell $ws.Range("E22") "  +1.52%  "
Set-TextCell $ws.Range("D23") "69.98"
Set-TextCell $ws.Range("E23") "  +0.06%  "
Set-TextCell $ws.Range("D24") "269.48"
Set-TextCell $ws.Range("E24") "  -2.38%  "
Set-TextCell $ws.Range("D25") "2.81"
Set-TextCell $ws.Range("E25") "  +2.07%  "
Set-TextCell $ws.Range("D26") "27.00"
Set-TextCell $ws.Range("E26") "  +1.28%  "
Set-TextCell $ws.Range("E27") "  +0.08%  "
Set-TextCell $ws.Range("D28") "10.31"
Set-TextCell $ws.Range("E28") "  +1.55%  "
Set-TextCell $ws.Range("E29") "  +1.67%  "
Set-TextCell $ws.Range("D33") "0.0470"
Set-TextCell $ws.Range("E33") "  +23.00%  "
Set-TextCell $ws.Range("D34") "5.89"
Set-TextCell $ws.Range("E34") "  +2.99%  "
Set-TextCell $ws.Range("D35") "5.44"
Set-TextCell $ws.Range("E35") "  +9.80%  "
Set-TextCell $ws.Range("D36") "0.0845"
Set-TextCell $ws.Range("E36") "  +3.44%  "
Set-TextCell $ws.Range("E37") "  +0.01%  "
Set-TextCell $ws.Range("D38") "3.28"
Set-TextCell $ws.Range("E38") "  +2.25%  "
Set-TextCell $ws.Range("D41") "0.116"
Set-TextCell $ws.Range("E41") "  +0.81%  "
Set-TextCell $ws.Range("D42") "23.43"
Set-TextCell $ws.Range("E42") "  +1.33%  "
Set-TextCell $ws.Range("E43") "  -4.22%  "
Set-TextCell $ws.Range("D49") "5.97"
Set-TextCell $ws.Range("E49") "  +7.57%  "
Set-TextCell $ws.Range("D50") "0.979"
Set-TextCell $ws.Range("E50") "  +11.00%  "
Set-TextCell $ws.Range("D51") "61.49"
Set-TextCell $ws.Range("E51") "  +3.57%  "

# --- Reordered block: rows 30-32 ---
Set-TextCell $ws.Range("B30") "InjectiveProtocol"
Set-TextCell $ws.Range("C30") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws.Range("D30") "35.57"
Set-TextCell $ws.Range("E30") "  +2.83%  "
Set-TextCell $ws.Range("B31") "OKB"
Set-TextCell $ws.Range("C31") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws.Range("D31") "52.93"
Set-TextCell $ws.Range("E31") "  +5.25%  "
Set-TextCell $ws.Range("B32") "Kaspa"
Set-TextCell $ws.Range("C32") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D32") "0.140"
Set-TextCell $ws.Range("E32") "  -0.80%  "

# --- Reordered block: rows 39-40 ---
Set-TextCell $ws.Range("B39") "Celestia"
Set-TextCell $ws.Range("C39") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell $ws.Range("D39") "18.61"
Set-TextCell $ws.Range("E39") "  -1.77%  "
Set-TextCell $ws.Range("B40") "ARBITRUM"
Set-TextCell $ws.Range("C40") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws.Range("D40") "2.05"
Set-TextCell $ws.Range("E40") "  -2.33%  "

# --- Reordered block: rows 44-47 ---
Set-TextCell $ws.Range("B44") "WEMIXToken"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws.Range("D44") "2.28"
Set-TextCell $ws.Range("E44") "  -2.97%  "
Set-TextCell $ws.Range("B45") "Monero"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D45") "123.93"
Set-TextCell $ws.Range("E45") "  -1.60%  "
Set-TextCell $ws.Range("B46") "Maker"
Set-TextCell $ws.Range("C46") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Range("D46") "2.099.36"
Set-TextCell $ws.Range("E46") "  +1.51%  "
Set-TextCell $ws.Range("B47") "NEARProtocol"
Set-TextCell $ws.Range("C47") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D47") "3.37"
Set-TextCell $ws.Range("E47") "  +1.88%  "
